$d = $word.ActiveDocument

# The report used to read:
#   "...После этого, была считана вся информация с сайта и попутно
#    переведена в lxml формат с помощью пакета анализа документов
#    BeautifulSoup"
# It should read just:
#   "...После этого, была считана вся информация с сайта"
# i.e. everything from "с сайта" (inclusive) through the end of
# "BeautifulSoup" is replaced by the plain phrase "с сайта".
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Text = "с сайта и попутно переведена в lxml формат с помощью пакета анализа документов BeautifulSoup"
$found = $r.Find.Execute()

if ($found) {
    # Delete the whole stretch, then retype the short replacement at the
    # same spot. Re-typing (rather than letting Find/Replace rewrite the
    # text in place) keeps the freshly inserted "с сайта" in its own run
    # instead of silently re-merging into the previous sentence's run,
    # matching the run layout that results from a real edit in Word.
    $r.Delete()
    $r.InsertAfter("с сайта")
    $r.Bold = $true
    $r.Bold = $false
    Write-Output "replaced trailing description with 'с сайта'"
} else {
    Write-Output "WARNING: target phrase not found, no changes made"
}
